$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as exact text
# (Excel would otherwise coerce them to floating-point numbers and lose formatting,
# e.g. "1.001" -> 1.0009999999999999). Temporarily force Text format, then restore
# the default "Normal" style so no stray per-cell style index is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.008.23'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.875.55'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '305.97'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.5066'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').Value = '0.3664'
$ws.Range('E8').Value = '  -2.05%  '
$ws.Range('D9').Value = '0.07203'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').Value = '0.8954'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').Value = '20.76'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D12').Value = '1.870.27'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = '0.07526'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').Value = '94.89'
$ws.Range('E14').Value = '  +6.51%  '
$ws.Range('D15').Value = '5.245'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '0.000008553'
$ws.Range('D18').Value = '14.28'
$ws.Range('E18').Value = '  +1.33%  '
$ws.Range('D19').Value = '0.9995'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '27.050.68'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = '5.035'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '2.112.90'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('D24').Value = '6.425'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '148.37'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').Value = '1.791'
$ws.Range('E26').Value = '  -2.68%  '
$ws.Range('D27').Value = '17.93'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = '2.082'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').Value = '113.46'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '4.710'
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('D31').Value = '4.691'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '0.09177'
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('D33').Value = '0.05154'
$ws.Range('E33').Value = '  +0.47%  '
$ws.Range('D34').Value = '0.7542'
$ws.Range('E34').Value = '  +4.15%  '
$ws.Range('D35').Value = '2.983'
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('D36').Value = '1.162'
$ws.Range('E36').Value = '  +1.05%  '
$ws.Range('E37').Value = '  +6.36%  '
$ws.Range('D38').Value = '2.584'
$ws.Range('E38').Value = '  +5.07%  '
$ws.Range('D39').Value = '0.5624'
$ws.Range('E39').Value = '  +6.65%  '
$ws.Range('E40').Value = '  -1.71%  '
$ws.Range('D41').Value = '1.073'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('D42').Value = '6.610'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').Value = '115.97'
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('D44').Value = '8.556'
$ws.Range('E44').Value = '  +3.80%  '
$ws.Range('D45').Value = '0.1479'
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('D46').Value = '0.4734'
$ws.Range('E46').Value = '  +2.89%  '
$ws.Range('D47').Value = '0.9994'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '10.11'
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('D49').Value = '1.567'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').Value = '36.91'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('D51').Value = '63.26'
$ws.Range('E51').Value = '  -0.96%  '

$ws.Range("D2:D51").Style = "Normal"

